# Add files via upload
#
# Adds two new survey-response columns (O, P) with Arabic headers matching
# the existing header style, and applies a custom timestamp number format
# to the first data row under the Timestamp column (A2), leaving the cell
# itself empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns O1 / P1 -------------------------------------------
# Clone the formatting of the last existing header cell (N1: bold font,
# thin border, centered alignment) onto the two new header cells, then set
# their text.
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)

$ws.Range("O1").Value = "الشعور تجاه التعليم الإلكتروني"
$ws.Range("P1").Value = "أسباب الإحباط"

# --- Timestamp-style number format on A2 -----------------------------------
# Give A2 a custom "yyyy-mm-dd hh:mm:ss" number format (matching the
# Timestamp column header above it) while leaving the cell value empty.
$ws.Range("A2").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("A2").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# --- Selection -------------------------------------------------------------
$ws.Range("D8").Select()
